# Update automàtic: dades i banners [2026-02-23 19:20]
# Refresh of the daily Meteocat extraction sheet: new DATA_EXTRACCIO
# timestamps plus the handful of measurements that shifted between runs.
# NOTE: values like "67%" are written with a leading '' (single-quote)
# prefix so Excel stores them as literal text, matching the original
# inlineStr cells, instead of auto-converting them to a percentage number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-23 19:18:55'
$ws.Range('O2').Value = '6.1 °C'
$ws.Range('E3').Value = '2026-02-23 19:18:57'
$ws.Range('E4').Value = '2026-02-23 19:19:00'
$ws.Range('H4').Value = '''67%'
$ws.Range('O4').Value = '12.4 °C'
$ws.Range('E5').Value = '2026-02-23 19:19:03'
$ws.Range('H5').Value = '''30%'
$ws.Range('E6').Value = '2026-02-23 19:19:05'
$ws.Range('H6').Value = '''60%'
$ws.Range('E7').Value = '2026-02-23 19:19:08'
$ws.Range('H7').Value = '''66%'
$ws.Range('E8').Value = '2026-02-23 19:19:11'
$ws.Range('H8').Value = '''57%'
$ws.Range('E9').Value = '2026-02-23 19:19:14'
$ws.Range('O9').Value = '12.7 °C'
$ws.Range('E10').Value = '2026-02-23 19:19:16'
$ws.Range('H10').Value = '''75%'
$ws.Range('O10').Value = '11.2 °C'
$ws.Range('E11').Value = '2026-02-23 19:19:17'
$ws.Range('E12').Value = '2026-02-23 19:19:19'
$ws.Range('E13').Value = '2026-02-23 19:19:20'
$ws.Range('E14').Value = '2026-02-23 19:19:21'
$ws.Range('H14').Value = '''73%'
$ws.Range('O14').Value = '13.0 °C'
$ws.Range('E15').Value = '2026-02-23 19:19:22'
$ws.Range('H15').Value = '''69%'
$ws.Range('E16').Value = '2026-02-23 19:19:23'
$ws.Range('H16').Value = '''20%'
$ws.Range('E17').Value = '2026-02-23 19:19:24'
$ws.Range('O17').Value = '8.7 °C'
$ws.Range('E18').Value = '2026-02-23 19:19:26'
$ws.Range('E19').Value = '2026-02-23 19:19:27'
$ws.Range('E20').Value = '2026-02-23 19:19:28'
$ws.Range('O20').Value = '4.4 °C'
$ws.Range('E21').Value = '2026-02-23 19:19:29'
$ws.Range('J21').Value = '1025.7 hPa'
$ws.Range('O21').Value = '9.8 °C'
$ws.Range('E22').Value = '2026-02-23 19:19:32'
$ws.Range('G22').Value = '109 cm'
$ws.Range('E23').Value = '2026-02-23 19:19:34'
$ws.Range('H23').Value = '''21%'
$ws.Range('E24').Value = '2026-02-23 19:19:37'
$ws.Range('O24').Value = '8.8 °C'
$ws.Range('E25').Value = '2026-02-23 19:19:40'
$ws.Range('O25').Value = '6.0 °C'
$ws.Range('E26').Value = '2026-02-23 19:19:42'
$ws.Range('G26').Value = '1 cm'
$ws.Range('O26').Value = '10.3 °C'
$ws.Range('E27').Value = '2026-02-23 19:19:44'
$ws.Range('H27').Value = '''28%'
$ws.Range('E28').Value = '2026-02-23 19:19:47'
$ws.Range('E29').Value = '2026-02-23 19:19:49'
$ws.Range('E30').Value = '2026-02-23 19:19:52'
$ws.Range('O30').Value = '13.2 °C'
$ws.Range('E31').Value = '2026-02-23 19:19:54'
$ws.Range('O31').Value = '16.5 °C'
$ws.Range('E32').Value = '2026-02-23 19:19:57'
$ws.Range('O32').Value = '8.3 °C'
$ws.Range('E33').Value = '2026-02-23 19:20:00'
$ws.Range('E34').Value = '2026-02-23 19:20:03'
$ws.Range('H34').Value = '''41%'
$ws.Range('O34').Value = '4.3 °C'
$ws.Range('E35').Value = '2026-02-23 19:20:05'
$ws.Range('K35').Value = '16.8 MJ/m2'
$ws.Range('O35').Value = '12.7 °C'
$ws.Range('E36').Value = '2026-02-23 19:20:08'
$ws.Range('E37').Value = '2026-02-23 19:20:11'
$ws.Range('E38').Value = '2026-02-23 19:20:13'
$ws.Range('H38').Value = '''63%'
$ws.Range('O38').Value = '12.4 °C'
$ws.Range('E39').Value = '2026-02-23 19:20:16'
$ws.Range('O39').Value = '4.2 °C'
$ws.Range('E40').Value = '2026-02-23 19:20:18'
$ws.Range('J40').Value = '1026.1 hPa'
$ws.Range('E41').Value = '2026-02-23 19:20:21'
$ws.Range('O41').Value = '12.3 °C'
$ws.Range('E42').Value = '2026-02-23 19:20:24'
$ws.Range('E43').Value = '2026-02-23 19:20:26'
$ws.Range('O43').Value = '10.5 °C'
$ws.Range('E44').Value = '2026-02-23 19:20:29'
$ws.Range('K44').Value = '16.2 MJ/m2'
$ws.Range('N44').Value = '0.9 °C 18:39 TU'
$ws.Range('E45').Value = '2026-02-23 19:20:32'
$ws.Range('K45').Value = '14.4 MJ/m2'
$ws.Range('O45').Value = '8.6 °C'
$ws.Range('E46').Value = '2026-02-23 19:20:34'
